$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "97÷3=32, 1"
$t.Cell(1, 2).Range.Text = "13÷8=1, 5"
$t.Cell(1, 3).Range.Text = "61÷5=12, 1"
$t.Cell(1, 4).Range.Text = "89÷6=14, 5"
$t.Cell(1, 5).Range.Text = "29÷3=9, 2"

$t.Cell(5, 1).Range.Text = "16÷3=5, 1"
$t.Cell(5, 2).Range.Text = "46÷7=6, 4"
$t.Cell(5, 3).Range.Text = "44÷7=6, 2"
$t.Cell(5, 4).Range.Text = "36÷4=9, 0"
$t.Cell(5, 5).Range.Text = "16÷2=8, 0"

$t.Cell(9, 1).Range.Text = "39÷6=6, 3"
$t.Cell(9, 2).Range.Text = "33÷4=8, 1"
$t.Cell(9, 3).Range.Text = "78÷7=11, 1"
$t.Cell(9, 4).Range.Text = "10÷9=1, 1"
$t.Cell(9, 5).Range.Text = "64÷2=32, 0"

$t.Cell(13, 1).Range.Text = "95÷9=10, 5"
$t.Cell(13, 2).Range.Text = "13÷8=1, 5"
$t.Cell(13, 3).Range.Text = "80÷6=13, 2"
$t.Cell(13, 4).Range.Text = "29÷5=5, 4"
$t.Cell(13, 5).Range.Text = "80÷5=16, 0"

$t.Cell(17, 1).Range.Text = "49÷8=6, 1"
$t.Cell(17, 2).Range.Text = "21÷5=4, 1"
$t.Cell(17, 3).Range.Text = "67÷4=16, 3"
$t.Cell(17, 4).Range.Text = "69÷2=34, 1"
$t.Cell(17, 5).Range.Text = "16÷3=5, 1"

